# Update NATMI LR-pair TPM-derived metrics (Uts2-Uts2r) with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.3735556666666667
$ws.Range("H2").Value = 1.120667
$ws.Range("I2").Value = 0.2390338626701809
$ws.Range("J2").Value = 0.239033862670181
$ws.Range("M2").Value = 0.04057
$ws.Range("O2").Value = 0.1939988045427376
$ws.Range("P2").Value = 0.1939988045427376
$ws.Range("Q2").Value = 0.01515515339666667
$ws.Range("R2").Value = 0.13639638057
$ws.Range("S2").Value = 0.04637228360324801
$ws.Range("T2").Value = 0.04637228360324801
$ws.Range("G3").Value = 0.3735556666666667
$ws.Range("H3").Value = 1.120667
$ws.Range("I3").Value = 0.2390338626701809
$ws.Range("J3").Value = 0.239033862670181
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.168555
$ws.Range("N3").Value = 0.505665
$ws.Range("O3").Value = 0.8060011954572625
$ws.Range("P3").Value = 0.8060011954572625
$ws.Range("Q3").Value = 0.062964675395
$ws.Range("R3").Value = 0.5666820785550001
$ws.Range("S3").Value = 0.1926615790669329
$ws.Range("T3").Value = 0.192661579066933
$ws.Range("G4").Value = 0.8080600000000001
$ws.Range("I4").Value = 0.5170680578689292
$ws.Range("J4").Value = 0.5170680578689292
$ws.Range("M4").Value = 0.04057
$ws.Range("O4").Value = 0.1939988045427376
$ws.Range("P4").Value = 0.1939988045427376
$ws.Range("Q4").Value = 0.03278299420000001
$ws.Range("S4").Value = 0.1003105850938073
$ws.Range("T4").Value = 0.1003105850938073
$ws.Range("G5").Value = 0.8080600000000001
$ws.Range("I5").Value = 0.5170680578689292
$ws.Range("J5").Value = 0.5170680578689292
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.168555
$ws.Range("N5").Value = 0.505665
$ws.Range("O5").Value = 0.8060011954572625
$ws.Range("P5").Value = 0.8060011954572625
$ws.Range("Q5").Value = 0.1362025533
$ws.Range("R5").Value = 1.2258229797
$ws.Range("S5").Value = 0.4167574727751219
$ws.Range("T5").Value = 0.4167574727751219
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3811573333333333
$ws.Range("H6").Value = 1.143472
$ws.Range("I6").Value = 0.2438980794608899
$ws.Range("J6").Value = 0.2438980794608899
$ws.Range("M6").Value = 0.04057
$ws.Range("O6").Value = 0.1939988045427376
$ws.Range("P6").Value = 0.1939988045427376
$ws.Range("Q6").Value = 0.01546355301333333
$ws.Range("R6").Value = 0.13917197712
$ws.Range("S6").Value = 0.04731593584568226
$ws.Range("T6").Value = 0.04731593584568227
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3811573333333333
$ws.Range("H7").Value = 1.143472
$ws.Range("I7").Value = 0.2438980794608899
$ws.Range("J7").Value = 0.2438980794608899
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.168555
$ws.Range("N7").Value = 0.505665
$ws.Range("O7").Value = 0.8060011954572625
$ws.Range("P7").Value = 0.8060011954572625
$ws.Range("Q7").Value = 0.06424597432000001
$ws.Range("R7").Value = 0.5782137688800001
$ws.Range("S7").Value = 0.1965821436152076
$ws.Range("T7").Value = 0.1965821436152077
